# "Generate Report for Handoff"
# Replaces the two pending localization rows (681b7454... / 890a3397...) with a
# fresh pair of handoff rows (8e10f827... / ffffc98e1ce6...) across the
# Overview / zh-cn / de-de sheets, updates the Latest Handoff / Handback
# bookkeeping columns, and drops the now-empty "Latest Target File" hyperlink
# column on the language sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---- capture the existing (unchanged) hyperlink target addresses before we
# ---- tear the hyperlink collections down and rebuild them. ----
$ov_b2_addr = $ws1.Range("B2").Hyperlinks.Item(1).Address
$ov_b3_addr = $ws1.Range("B3").Hyperlinks.Item(1).Address

$zh_a2_addr = $ws2.Range("A2").Hyperlinks.Item(1).Address
$zh_a3_addr = $ws2.Range("A3").Hyperlinks.Item(1).Address

$de_a2_addr = $ws3.Range("A2").Hyperlinks.Item(1).Address
$de_a3_addr = $ws3.Range("A3").Hyperlinks.Item(1).Address

# New identifiers generated for this handoff.
$newFile1 = "8e10f827-73d5-486c-8e88-90377d54518d.md"
$newFile2 = "ffffc98e1ce6-4692-4345-bd2e-0288eee64758.md"

$newXlfZh = "8e10f827-73d5-486c-8e88-90377d54518d.dbf721f12bda5925636c4c25f580973fc6fd0f50.zh-cn.xlf"
$newXlfDe = "8e10f827-73d5-486c-8e88-90377d54518d.dbf721f12bda5925636c4c25f580973fc6fd0f50.de-de.xlf"

$status = "Ready for handoff"
$genDate = "2016-08-27 01:02:10"
$handoffDateZh = "2016-08-27 01:01:59"
$handbackDate = "0001-01-01 00:00:00"

# ================= Overview sheet =================
$ws1.Range("A2").Value2 = $newFile1
$ws1.Range("B2").Value2 = "e2e\" + $newFile1
$ws1.Range("E2").Value2 = $status
$ws1.Range("F2").Value2 = $status
$ws1.Range("G2").Value2 = $genDate

$ws1.Range("A3").Value2 = $newFile2
$ws1.Range("B3").Value2 = "e2e\" + $newFile2
$ws1.Range("E3").Value2 = $status
$ws1.Range("F3").Value2 = $status
$ws1.Range("G3").Value2 = $genDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $ov_b2_addr, "", "", "e2e\" + $newFile1)
$ws1.Hyperlinks.Add($ws1.Range("B3"), $ov_b3_addr, "", "", "e2e\" + $newFile2)

# ================= zh-cn sheet =================
$ws2.Range("A2").Value2 = $newFile1
$ws2.Range("C2").Value2 = $status
$ws2.Range("G2").Value2 = $newXlfZh
$ws2.Range("H2").Value2 = $handoffDateZh
$ws2.Range("I2").Value2 = ""
$ws2.Range("J2").Value2 = ""
$ws2.Range("K2").Value2 = $handbackDate

$ws2.Range("A3").Value2 = $newFile2
$ws2.Range("C3").Value2 = $status
$ws2.Range("F3").Value2 = "True"
$ws2.Range("G3").Value2 = $newXlfZh
$ws2.Range("H3").Value2 = $handoffDateZh
$ws2.Range("I3").Value2 = ""
$ws2.Range("J3").Value2 = ""
$ws2.Range("K3").Value2 = $handbackDate

$ws2.Range("I2").Style = "Normal"
$ws2.Range("I3").Style = "Normal"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $zh_a2_addr, "", "", $newFile1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $zh_a3_addr, "", "", $newFile2)

# ================= de-de sheet =================
$ws3.Range("A2").Value2 = $newFile1
$ws3.Range("C2").Value2 = $status
$ws3.Range("G2").Value2 = $newXlfDe
$ws3.Range("H2").Value2 = $genDate
$ws3.Range("I2").Value2 = ""
$ws3.Range("J2").Value2 = ""
$ws3.Range("K2").Value2 = $handbackDate

$ws3.Range("A3").Value2 = $newFile2
$ws3.Range("C3").Value2 = $status
$ws3.Range("F3").Value2 = "True"
$ws3.Range("G3").Value2 = $newXlfDe
$ws3.Range("H3").Value2 = $genDate
$ws3.Range("I3").Value2 = ""
$ws3.Range("J3").Value2 = ""
$ws3.Range("K3").Value2 = $handbackDate

$ws3.Range("I2").Style = "Normal"
$ws3.Range("I3").Style = "Normal"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $de_a2_addr, "", "", $newFile1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $de_a3_addr, "", "", $newFile2)

# ================= column widths (cosmetic re-fit after content change) =================
$ws1.Columns.Item(5).ColumnWidth = 17.2159881591797
$ws1.Columns.Item(6).ColumnWidth = 17.2159881591797

$ws2.Columns.Item(3).ColumnWidth = 17.2159881591797
$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426

$ws3.Columns.Item(3).ColumnWidth = 17.2159881591797
$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
